# Replace the hard-coded "2022. Quincentenario..." legend in the letterhead
# paragraph with the ${leyenda} placeholder.
#
# In the source document this sentence is split across three runs
# (with <w:proofErr/> spell-check markers sitting between them):
#   run1: "\u201c2022. A\u00f1o del "
#   run2: "Quincentenario"
#   run3: " de Toluca, Capital del Estado de M\u00e9xico\"."
# All three runs share identical rPr (Gotham HTF, color 808080, sz 18, szCs 20).
#
# Find.Execute can match text that spans multiple runs/fields inside a
# paragraph, and replacing it collapses the match into a single run that
# inherits the formatting of the first run, while the proofErr markers
# (which carry no text) are dropped along with the now-empty extra runs.

$d = $word.ActiveDocument

# Build the search string without relying on the source file's literal
# encoding: U+201C (left double quotation mark), "n with tilde" (U+00F1),
# "e with acute" (U+00E9), and a plain ASCII double quote before the
# trailing period.
$oldText = [string]::Concat(
    [char]0x201C,
    "2022. A",
    [char]0x00F1,
    "o del Quincentenario de Toluca, Capital del Estado de M",
    [char]0x00E9,
    "xico",
    [char]0x0022,
    "."
)

$newText = '${leyenda}'

$found = $d.Content.Find.Execute(
    $oldText,   # FindText
    $true,      # MatchCase
    $false,     # MatchWholeWord
    $false,     # MatchWildcards
    $false,     # MatchSoundsLike
    $false,     # MatchAllWordForms
    $true,      # Forward
    1,          # Wrap (wdFindContinue)
    $false,     # Format
    $newText,   # ReplaceWith
    2           # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Could not find the legend text to replace."
}
